$d = $word.ActiveDocument

# Locate the "Author" styled paragraph whose text is exactly "Edison Achalma"
# (the author byline right under the document title) so we can add a new
# affiliation paragraph directly after it.
$targetEnd = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author" -and $p.Range.Text -eq "Edison Achalma`r") {
        $targetEnd = $p.Range.End
    }
}

if ($targetEnd -ge 0) {
    # Build a brand-new, detached range collapsed right after the existing
    # text (i.e. right before its paragraph mark) and split a new paragraph
    # there. This leaves the existing "Edison Achalma" run/paragraph
    # completely untouched.
    $insertionPoint = $d.Range($targetEnd, $targetEnd)
    $insertionPoint.InsertParagraphAfter()

    # Re-fetch the freshly created (still empty) paragraph and populate it.
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -eq $targetEnd) {
            $p.Range.InsertAfter("Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga")
            $p.Style = "Author"
            break
        }
    }
}
